# Update res_bus/vm_pu.xlsx results for "case with 380 kV done"
# Slack-bus voltage setpoint (B column) changed from 1.05 pu to 1.02 pu,
# with the downstream bus voltage-magnitude results (C:F, I:N) recomputed accordingly
# for every row (r=2..25). Columns G (=1) and H (blank) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.033566086006686
$ws.Range("D2").Value = 1.041284775085484
$ws.Range("E2").Value = 1.037209081715803
$ws.Range("F2").Value = 1.04844955719123
$ws.Range("I2").Value = 1.039348246348588
$ws.Range("J2").Value = 1.03869011018239
$ws.Range("K2").Value = 1.044064724237527
$ws.Range("L2").Value = 1.040000619992905
$ws.Range("M2").Value = 1.051209370383304
$ws.Range("N2").Value = 1.016864039476882

# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.034406165869343
$ws.Range("D3").Value = 1.041946689113928
$ws.Range("E3").Value = 1.037997028747212
$ws.Range("F3").Value = 1.049291378171593
$ws.Range("I3").Value = 1.039564980280528
$ws.Range("J3").Value = 1.039173511798313
$ws.Range("K3").Value = 1.044537625248489
$ws.Range("L3").Value = 1.040598385274516
$ws.Range("M3").Value = 1.05186316195899
$ws.Range("N3").Value = 1.017025386728734

# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.034950295680202
$ws.Range("D4").Value = 1.042375425894446
$ws.Range("E4").Value = 1.038507759910976
$ws.Range("F4").Value = 1.049837006658326
$ws.Range("I4").Value = 1.039704127733229
$ws.Range("J4").Value = 1.039486187120999
$ws.Range("K4").Value = 1.044843373994666
$ws.Range("L4").Value = 1.040985399239814
$ws.Range("M4").Value = 1.052286475353216
$ws.Range("N4").Value = 1.017129715324971

# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.03517917604336
$ws.Range("D5").Value = 1.042555769047492
$ws.Range("E5").Value = 1.038722679397158
$ws.Range("F5").Value = 1.050066605824337
$ws.Range("I5").Value = 1.039762362655128
$ws.Range("J5").Value = 1.039617606398579
$ws.Range("K5").Value = 1.04497184968291
$ws.Range("L5").Value = 1.041148151102138
$ws.Range("M5").Value = 1.052464498539729
$ws.Range("N5").Value = 1.017173556935251

# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.035217613535415
$ws.Range("D6").Value = 1.04258605539994
$ws.Range("E6").Value = 1.038758777475903
$ws.Range("F6").Value = 1.050105169202333
$ws.Range("I6").Value = 1.03977212512898
$ws.Range("J6").Value = 1.039639670512551
$ws.Range("K6").Value = 1.044993417697184
$ws.Range("L6").Value = 1.041175480815897
$ws.Range("M6").Value = 1.052494393019916
$ws.Range("N6").Value = 1.017180917058248

# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.034953353489058
$ws.Range("D7").Value = 1.042377835249033
$ws.Range("E7").Value = 1.038510630860205
$ws.Range("F7").Value = 1.049840073723496
$ws.Range("I7").Value = 1.039704906904188
$ws.Range("J7").Value = 1.039487943268972
$ws.Range("K7").Value = 1.044845090934719
$ws.Range("L7").Value = 1.040987573738826
$ws.Range("M7").Value = 1.052288853865173
$ws.Range("N7").Value = 1.01713030121052

# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.033849881775224
$ws.Range("D8").Value = 1.041508381102153
$ws.Range("E8").Value = 1.037475189492297
$ws.Range("F8").Value = 1.048733864155479
$ws.Range("I8").Value = 1.039421718620189
$ws.Range("J8").Value = 1.038853501692106
$ws.Range("K8").Value = 1.044224594310834
$ws.Range("L8").Value = 1.040202591178395
$ws.Range("M8").Value = 1.051430266176398
$ws.Range("N8").Value = 1.016918582560855

# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.031909643321652
$ws.Range("D9").Value = 1.039979692996363
$ws.Range("E9").Value = 1.035657403447846
$ws.Range("F9").Value = 1.046791664475473
$ws.Range("I9").Value = 1.038914361571618
$ws.Range("J9").Value = 1.037734685721066
$ws.Range("K9").Value = 1.04312934292384
$ws.Range("L9").Value = 1.038821099224714
$ws.Range("M9").Value = 1.049919435715883
$ws.Range("N9").Value = 1.016544961335514

# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.030619080615271
$ws.Range("D10").Value = 1.038962957038768
$ws.Range("E10").Value = 1.034450220555741
$ws.Range("F10").Value = 1.045501738128985
$ws.Range("I10").Value = 1.038570560931964
$ws.Range("J10").Value = 1.036988309979252
$ws.Range("K10").Value = 1.042398002612382
$ws.Range("L10").Value = 1.037901365463128
$ws.Range("M10").Value = 1.048913728114078
$ws.Range("N10").Value = 1.016295540141535

# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030060966967851
$ws.Range("D11").Value = 1.038523287547333
$ws.Range("E11").Value = 1.033928628234081
$ws.Range("F11").Value = 1.044944365148994
$ws.Range("I11").Value = 1.038420381909431
$ws.Range("J11").Value = 1.036665018025183
$ws.Range("K11").Value = 1.042081061992501
$ws.Range("L11").Value = 1.037503427733365
$ws.Range("M11").Value = 1.048478624114106
$ws.Range("N11").Value = 1.016187462890646

# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.029853766703696
$ws.Range("D12").Value = 1.038360064193133
$ws.Range("E12").Value = 1.033735056582497
$ws.Range("F12").Value = 1.044737510227136
$ws.Range("I12").Value = 1.038364402530334
$ws.Range("J12").Value = 1.036544918239391
$ws.Range("K12").Value = 1.041963297532064
$ws.Range("L12").Value = 1.037355664428431
$ws.Range("M12").Value = 1.048317064892465
$ws.Range("N12").Value = 1.016147307184682

# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.029898206951863
$ws.Range("D13").Value = 1.038395072079869
$ws.Range("E13").Value = 1.033776570593376
$ws.Range("F13").Value = 1.04478187321483
$ws.Range("I13").Value = 1.038376419163398
$ws.Range("J13").Value = 1.036570680699507
$ws.Range("K13").Value = 1.041988560134932
$ws.Range("L13").Value = 1.037387357941516
$ws.Range("M13").Value = 1.048351717253982
$ws.Range("N13").Value = 1.016155921210611

# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030043837517364
$ws.Range("D14").Value = 1.038509793618666
$ws.Range("E14").Value = 1.033912624042931
$ws.Range("F14").Value = 1.044927262810041
$ws.Range("I14").Value = 1.038415758634516
$ws.Range("J14").Value = 1.036655090833481
$ws.Range("K14").Value = 1.042071328325559
$ws.Range("L14").Value = 1.037491212574983
$ws.Range("M14").Value = 1.048465268390657
$ws.Range("N14").Value = 1.01618414382932

# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.030133579671842
$ws.Range("D15").Value = 1.038580489253078
$ws.Range("E15").Value = 1.033996473771781
$ws.Range("F15").Value = 1.045016865810622
$ws.Range("I15").Value = 1.038439971005275
$ws.Range("J15").Value = 1.036707096783503
$ws.Range("K15").Value = 1.042122319462489
$ws.Range("L15").Value = 1.037555207310618
$ws.Range("M15").Value = 1.04853523869699
$ws.Range("N15").Value = 1.016201531271628

# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03065613583967
$ws.Range("D16").Value = 1.038992148906058
$ws.Range("E16").Value = 1.034484860795456
$ws.Range("F16").Value = 1.045538754026915
$ws.Range("I16").Value = 1.038580500261987
$ws.Range("J16").Value = 1.037009763664979
$ws.Range("K16").Value = 1.042419031434
$ws.Range("L16").Value = 1.037927781999454
$ws.Range("M16").Value = 1.048942612549086
$ws.Range("N16").Value = 1.016302711305077

# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.030984112300457
$ws.Range("D17").Value = 1.039250529688915
$ws.Range("E17").Value = 1.034791515773369
$ws.Range("F17").Value = 1.04586643627071
$ws.Range("I17").Value = 1.038668300083326
$ws.Range("J17").Value = 1.037199590993364
$ws.Range("K17").Value = 1.042605080850542
$ws.Range("L17").Value = 1.038161573098599
$ws.Range("M17").Value = 1.049198248659679
$ws.Range("N17").Value = 1.016366158775982

# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03117548367918
$ws.Range("D18").Value = 1.039401295047643
$ws.Range("E18").Value = 1.034970490938972
$ws.Range("F18").Value = 1.046057680891864
$ws.Range("I18").Value = 1.038719385681964
$ws.Range("J18").Value = 1.037310303721303
$ws.Range("K18").Value = 1.042713574630889
$ws.Range("L18").Value = 1.038297969599491
$ws.Range("M18").Value = 1.049347392860096
$ws.Range("N18").Value = 1.016403159227293

# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.03124074790457
$ws.Range("D19").Value = 1.03945271162142
$ws.Range("E19").Value = 1.035031535194397
$ws.Range("F19").Value = 1.046122909493661
$ws.Range("I19").Value = 1.038736783062131
$ws.Range("J19").Value = 1.037348052109443
$ws.Range("K19").Value = 1.042750563808819
$ws.Range("L19").Value = 1.038344482344169
$ws.Range("M19").Value = 1.049398253233225
$ws.Range("N19").Value = 1.016415774158703

# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.030948916458728
$ws.Range("D20").Value = 1.039222802048534
$ws.Range("E20").Value = 1.034758603360085
$ws.Range("F20").Value = 1.045831267353702
$ws.Range("I20").Value = 1.038658893077697
$ws.Range("J20").Value = 1.037179225378386
$ws.Range("K20").Value = 1.042585122157606
$ws.Range("L20").Value = 1.038136486401521
$ws.Range("M20").Value = 1.049170817591099
$ws.Range("N20").Value = 1.016359352223437

# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030000949951952
$ws.Range("D21").Value = 1.038476008500975
$ws.Range("E21").Value = 1.03387255495915
$ws.Range("F21").Value = 1.044884444263042
$ws.Range("I21").Value = 1.038404179550757
$ws.Range("J21").Value = 1.036630234550562
$ws.Range("K21").Value = 1.042046956211726
$ws.Range("L21").Value = 1.037460628623819
$ws.Range("M21").Value = 1.048431828803452
$ws.Range("N21").Value = 1.016175833271664

# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.029405551375574
$ws.Range("D22").Value = 1.038006988274285
$ws.Range("E22").Value = 1.033316451185944
$ws.Range("F22").Value = 1.044290171069428
$ws.Range("I22").Value = 1.038242896257252
$ws.Range("J22").Value = 1.036284977435884
$ws.Range("K22").Value = 1.041708367211929
$ws.Range("L22").Value = 1.037035971258617
$ws.Range("M22").Value = 1.047967532285625
$ws.Range("N22").Value = 1.016060384186074

# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.029721123609939
$ws.Range("D23").Value = 1.038255575037869
$ws.Range("E23").Value = 1.033611157868889
$ws.Range("F23").Value = 1.044605108073499
$ws.Range("I23").Value = 1.038328502889774
$ws.Range("J23").Value = 1.036468012396134
$ws.Range("K23").Value = 1.041887880331075
$ws.Range("L23").Value = 1.037261063018998
$ws.Range("M23").Value = 1.04821363236568
$ws.Range("N23").Value = 1.016121591793762

# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03096481973389
$ws.Range("D24").Value = 1.039235330802212
$ws.Range("E24").Value = 1.034773474727916
$ws.Range("F24").Value = 1.045847158324207
$ws.Range("I24").Value = 1.038663144089733
$ws.Range("J24").Value = 1.037188427755265
$ws.Range("K24").Value = 1.042594140711385
$ws.Range("L24").Value = 1.038147821907235
$ws.Range("M24").Value = 1.049183212399249
$ws.Range("N24").Value = 1.016362427834421

# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.032410731891767
$ws.Range("D25").Value = 1.040374481867992
$ws.Range("E25").Value = 1.036126528607396
$ws.Range("F25").Value = 1.047292918816677
$ws.Range("I25").Value = 1.039046509411485
$ws.Range("J25").Value = 1.038024019611997
$ws.Range("K25").Value = 1.043412703913697
$ws.Range("L25").Value = 1.039178032057877
$ws.Range("M25").Value = 1.050309761829035
$ws.Range("N25").Value = 1.016641613322807

Write-Output "Updated vm_pu.xlsx rows 2-25 (B,C,D,E,F,I,J,K,L,M,N) for 380 kV slack setpoint 1.02 pu"
